# Update "广州-漫展信息" workbook with freshly scraped counts (gh-pages output
# regenerated at commit 456a3b4).
#
# Sheet map:
#   1 "展览"     (Exhibitions)
#   2 "演出"     (Performances)
#   3 "本地生活" (Local life)
#   4 "全部类型" (All types - union of the three sheets above)
#
# Column F = "想去人数" (number of people interested) - numeric, grows over time.
# Column G = "最低票价" (minimum ticket price) - numeric, or text when the item
#            becomes unavailable for sale.

$wb = $excel.ActiveWorkbook

$sheetExhibition  = $wb.Worksheets.Item("展览")
$sheetPerformance = $wb.Worksheets.Item("演出")
$sheetLocalLife   = $wb.Worksheets.Item("本地生活")
$sheetAll         = $wb.Worksheets.Item("全部类型")

# ---- 展览 (sheet 1) ----
$sheetExhibition.Range("F4").Value  = 189
$sheetExhibition.Range("F5").Value  = 708
$sheetExhibition.Range("F6").Value  = 198
$sheetExhibition.Range("F7").Value  = 506
$sheetExhibition.Range("F8").Value  = 95
$sheetExhibition.Range("F9").Value  = 540
$sheetExhibition.Range("F10").Value = 461
$sheetExhibition.Range("F11").Value = 73
$sheetExhibition.Range("F12").Value = 34
$sheetExhibition.Range("F13").Value = 129
$sheetExhibition.Range("F14").Value = 209

# ---- 演出 (sheet 2) ----
$sheetPerformance.Range("F15").Value = 41
$sheetPerformance.Range("F16").Value = 26

# ---- 本地生活 (sheet 3) ----
$sheetLocalLife.Range("F2").Value = 6238
$sheetLocalLife.Range("G3").Value = "不可售"
$sheetLocalLife.Range("F4").Value = 763
$sheetLocalLife.Range("F5").Value = 1847

# ---- 全部类型 (sheet 4, mirrors the three sheets above) ----
$sheetAll.Range("F2").Value  = 6238
$sheetAll.Range("G3").Value  = "不可售"
$sheetAll.Range("F4").Value  = 763
$sheetAll.Range("F5").Value  = 1847
$sheetAll.Range("F11").Value = 189
$sheetAll.Range("F14").Value = 708
$sheetAll.Range("F16").Value = 198
$sheetAll.Range("F18").Value = 507
$sheetAll.Range("F20").Value = 95
$sheetAll.Range("F21").Value = 540
$sheetAll.Range("F23").Value = 461
$sheetAll.Range("F24").Value = 73
$sheetAll.Range("F27").Value = 34
$sheetAll.Range("F28").Value = 129
$sheetAll.Range("F31").Value = 41
$sheetAll.Range("F32").Value = 26
$sheetAll.Range("F34").Value = 209
